$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TSVData filename cell (D2): the export naming convention changed
# from a dot-separated suffix ("...fastq.TSVData.xlsx") to an
# underscore-separated suffix ("...fastq_TSVData.xlsx").
$ws.Range("D2").Value = "TC02_CDS_phs002430_Sex-NotSpecified_Experimental-strategy-RNA-Seq_Filetype-fastq_TSVData.xlsx"

# Reflect the updated view/selection state used when the workbook was saved.
$ws.Range("D2").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 2
